$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ARM")
# Row 88
$ws.Range("H88").Value = 2781.4
$ws.Range("I88").Value = 3000
$ws.Range("J88").Value = 2757.111
$ws.Range("K88").Value = 3000
$ws.Range("L88").Value = 2757.111
$ws.Range("M88").Value = -2594
$ws.Range("N88").Value = -3569.111

# Row 91
$ws.Range("H91").Value = 2781.4
$ws.Range("I91").Value = 3000
$ws.Range("J91").Value = 2757.111
$ws.Range("K91").Value = 3000
$ws.Range("L91").Value = 2757.111
$ws.Range("M91").Value = -1596
$ws.Range("N91").Value = -5565.111

# Row 122
$ws.Range("H122").Value = 10246.154
$ws.Range("I122").Value = 13630.223
$ws.Range("J122").Value = 2632
$ws.Range("K122").Value = 40890.669
$ws.Range("L122").Value = 7896
$ws.Range("M122").Value = -38440.669
$ws.Range("N122").Value = -12796

# Row 135
$ws.Range("H135").Value = 19817.125
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 19817.125
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 19817.125
$ws.Range("N135").Value = -29957.125

$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value = 33335172
$ws.Range("I86").Value = 1955.5555
$ws.Range("J86").Value = 83335000
$ws.Range("K86").Value = 1955.5555
$ws.Range("L86").Value = 83335000
$ws.Range("M86").Value = -832.5554999999999
$ws.Range("N86").Value = -83337246

# Row 89
$ws.Range("H89").Value = 33335172
$ws.Range("I89").Value = 1955.5555
$ws.Range("J89").Value = 83335000
$ws.Range("K89").Value = 9777.7775
$ws.Range("L89").Value = 416675000
$ws.Range("M89").Value = -4161.7775
$ws.Range("N89").Value = -416686232

$ws = $wb.Worksheets.Item("CRP")
# Row 62
$ws.Range("H62").Value = 2557
$ws.Range("I62").Value = 2354.5454
$ws.Range("J62").Value = 3002.4
$ws.Range("K62").Value = 2354.5454
$ws.Range("L62").Value = 3002.4
$ws.Range("M62").Value = -1730.5454
$ws.Range("N62").Value = -4250.4

# Row 65
$ws.Range("H65").Value = 2557
$ws.Range("I65").Value = 2354.5454
$ws.Range("J65").Value = 3002.4
$ws.Range("K65").Value = 11772.727
$ws.Range("L65").Value = 15012
$ws.Range("M65").Value = -8652.726999999999
$ws.Range("N65").Value = -21252

# Row 124
$ws.Range("H124").Value = 26244.5
$ws.Range("I124").Value = 0
$ws.Range("J124").Value = 26244.5
$ws.Range("K124").Value = 0
$ws.Range("L124").Value = 26244.5
$ws.Range("N124").Value = -31154.5

# Row 131
$ws.Range("H131").Value = 41000
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 41000
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 41000
$ws.Range("N131").Value = -51080

$ws = $wb.Worksheets.Item("GSM")
# Row 42
$ws.Range("H42").Value = 0
$ws.Range("I42").Value = 0
$ws.Range("J42").Value = 0
$ws.Range("K42").Value = 0
$ws.Range("L42").Value = 0
$ws.Range("N42").Value = ""

# Row 51
$ws.Range("H51").Value = 22940
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 22940
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 22940
$ws.Range("N51").Value = -23958

# Row 115
$ws.Range("H115").Value = 0
$ws.Range("I115").Value = 0
$ws.Range("J115").Value = 0
$ws.Range("K115").Value = 0
$ws.Range("L115").Value = 0
$ws.Range("N115").Value = ""

# Row 125
$ws.Range("H125").Value = 52244.5
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 52244.5
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 52244.5
$ws.Range("N125").Value = -57164.5

$ws = $wb.Worksheets.Item("LTW")
# Row 6
$ws.Range("H6").Value = 30241
$ws.Range("I6").Value = 0
$ws.Range("J6").Value = 30241
$ws.Range("K6").Value = 0
$ws.Range("L6").Value = 30241
$ws.Range("N6").Value = -30465

# Row 16
$ws.Range("H16").Value = 8500
$ws.Range("I16").Value = 5625
$ws.Range("J16").Value = 20000
$ws.Range("K16").Value = 5625
$ws.Range("L16").Value = 20000
$ws.Range("M16").Value = -5455
$ws.Range("N16").Value = -20340

# Row 22
$ws.Range("H22").Value = 1036.375
$ws.Range("I22").Value = 1058.2
$ws.Range("J22").Value = 1000
$ws.Range("K22").Value = 1058.2
$ws.Range("L22").Value = 1000
$ws.Range("M22").Value = -763.2
$ws.Range("N22").Value = -1590

# Row 27
$ws.Range("H27").Value = 1036.375
$ws.Range("I27").Value = 1058.2
$ws.Range("J27").Value = 1000
$ws.Range("K27").Value = 1058.2
$ws.Range("L27").Value = 1000
$ws.Range("M27").Value = -951.2
$ws.Range("N27").Value = -1214

# Row 46
$ws.Range("H46").Value = 554.381
$ws.Range("I46").Value = 614
$ws.Range("J46").Value = 500.18182
$ws.Range("K46").Value = 614
$ws.Range("L46").Value = 500.18182
$ws.Range("M46").Value = -426
$ws.Range("N46").Value = -876.18182

# Row 55
$ws.Range("H55").Value = 673
$ws.Range("I55").Value = 333.33334
$ws.Range("J55").Value = 876.8
$ws.Range("K55").Value = 333.33334
$ws.Range("L55").Value = 876.8
$ws.Range("M55").Value = -160.33334
$ws.Range("N55").Value = -1222.8

# Row 115
$ws.Range("H115").Value = 50302
$ws.Range("I115").Value = 0
$ws.Range("J115").Value = 50302
$ws.Range("K115").Value = 0
$ws.Range("L115").Value = 50302
$ws.Range("N115").Value = -52652

# Row 122
$ws.Range("H122").Value = 7504.8823
$ws.Range("I122").Value = 9870
$ws.Range("J122").Value = 5849.3
$ws.Range("K122").Value = 29610
$ws.Range("L122").Value = 17547.9
$ws.Range("M122").Value = -27160
$ws.Range("N122").Value = -22447.9

# Row 134
$ws.Range("H134").Value = 15571.077
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 15571.077
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 15571.077
$ws.Range("N134").Value = -25711.077

# Row 137
$ws.Range("H137").Value = 36809.668
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 36809.668
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 36809.668
$ws.Range("M137").Value = ""
$ws.Range("N137").Value = -47009.668

$ws = $wb.Worksheets.Item("WVR")
# Row 81
$ws.Range("H81").Value = 2552.3845
$ws.Range("I81").Value = 2084.875
$ws.Range("J81").Value = 3300.4
$ws.Range("K81").Value = 4169.75
$ws.Range("L81").Value = 6600.8
$ws.Range("M81").Value = -3108.75
$ws.Range("N81").Value = -8722.799999999999

# Row 84
$ws.Range("H84").Value = 2552.3845
$ws.Range("I84").Value = 2084.875
$ws.Range("J84").Value = 3300.4
$ws.Range("K84").Value = 20848.75
$ws.Range("L84").Value = 33004
$ws.Range("M84").Value = -15544.75
$ws.Range("N84").Value = -43612

# Row 122
$ws.Range("H122").Value = 2799.1428
$ws.Range("I122").Value = 2574
$ws.Range("J122").Value = 4150
$ws.Range("K122").Value = 7722
$ws.Range("L122").Value = 12450
$ws.Range("M122").Value = -5272
$ws.Range("N122").Value = -17350

# Row 126
$ws.Range("H126").Value = 4371.3335
$ws.Range("I126").Value = 2417.75
$ws.Range("J126").Value = 20000
$ws.Range("K126").Value = 7253.25
$ws.Range("L126").Value = 60000
$ws.Range("M126").Value = -4783.25
$ws.Range("N126").Value = -64940
